$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.919.72'
$ws.Range("E2").Value = '  +1.99%  '
$ws.Range("D3").Value = '1.811.69'
$ws.Range("E3").Value = '  +1.39%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").Value = '''337.82'
$ws.Range("E5").Value = '  +0.58%  '
$ws.Range("D6").Value = '''0.9966'
$ws.Range("E6").Value = '  -0.52%  '
$ws.Range("D7").Value = '''0.3917'
$ws.Range("E7").Value = '  +3.41%  '
$ws.Range("D8").Value = '''0.3479'
$ws.Range("E8").Value = '  +1.84%  '
$ws.Range("D9").Value = '''48.26'
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("D10").Value = '''1.200'
$ws.Range("E10").Value = '  +0.08%  '
$ws.Range("D11").Value = '''0.07560'
$ws.Range("E11").Value = '  +1.39%  '
$ws.Range("D12").Value = '''0.9968'
$ws.Range("E12").Value = '  -0.63%  '
$ws.Range("D13").Value = '''22.16'
$ws.Range("E13").Value = '  +1.00%  '
$ws.Range("D14").Value = '''6.515'
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("D15").Value = '1.812.52'
$ws.Range("E15").Value = '  +1.53%  '
$ws.Range("D16").Value = '''7.195'
$ws.Range("E16").Value = '  +2.58%  '
$ws.Range("D17").Value = '''0.00001107'
$ws.Range("E17").Value = '  +1.50%  '
$ws.Range("D18").Value = '''0.06670'
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("D19").Value = '''85.06'
$ws.Range("E19").Value = '  +0.96%  '
$ws.Range("D20").Value = '''0.9962'
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("D21").Value = '''17.86'
$ws.Range("E21").Value = '  +3.15%  '
$ws.Range("D22").Value = '''6.576'
$ws.Range("E22").Value = '  +1.82%  '
$ws.Range("D23").Value = '27.949.67'
$ws.Range("E23").Value = '  +2.25%  '
$ws.Range("D24").Value = '''12.87'
$ws.Range("E24").Value = '  +3.15%  '
$ws.Range("D25").Value = '''2.400'
$ws.Range("E25").Value = '  -2.09%  '
$ws.Range("D26").Value = '''2.551'
$ws.Range("E26").Value = '  +0.28%  '
$ws.Range("D27").Value = '''1.473'
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("D28").Value = '''21.32'
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("D29").Value = '''155.10'
$ws.Range("E29").Value = '  +3.27%  '
$ws.Range("D30").Value = '2.018.48'
$ws.Range("E30").Value = '  +1.61%  '
$ws.Range("D31").Value = '''135.57'
$ws.Range("E31").Value = '  +1.99%  '
$ws.Range("D32").Value = '''4.038'
$ws.Range("E32").Value = '  -0.70%  '
$ws.Range("D33").Value = '''6.121'
$ws.Range("E33").Value = '  +0.36%  '
$ws.Range("D34").Value = '''0.08831'
$ws.Range("E34").Value = '  +2.28%  '
$ws.Range("D35").Value = '''13.28'
$ws.Range("E35").Value = '  +0.75%  '
$ws.Range("D36").Value = '''5.551'
$ws.Range("E36").Value = '  +2.55%  '
$ws.Range("B37").Value = 'TheSandbox'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D37").Value = '''0.6935'
$ws.Range("E37").Value = '  +1.04%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.02421'
$ws.Range("E38").Value = '  +3.60%  '
$ws.Range("D39").Value = '''0.06520'
$ws.Range("E39").Value = '  +3.03%  '
$ws.Range("D40").Value = '''1.614'
$ws.Range("E40").Value = '  -3.45%  '
$ws.Range("D41").Value = '''0.2227'
$ws.Range("E41").Value = '  +1.76%  '
$ws.Range("D42").Value = '''1.266'
$ws.Range("E42").Value = '  -0.69%  '
$ws.Range("D43").Value = '''8.543'
$ws.Range("E43").Value = '  -2.94%  '
$ws.Range("D44").Value = '''14.83'
$ws.Range("E44").Value = '  +3.43%  '
$ws.Range("D45").Value = '''0.6530'
$ws.Range("E45").Value = '  +1.78%  '
$ws.Range("D46").Value = '''0.9955'
$ws.Range("E46").Value = '  -0.68%  '
$ws.Range("D47").Value = '''3.868'
$ws.Range("E47").Value = '  +0.54%  '
$ws.Range("D48").Value = '''2.164'
$ws.Range("E48").Value = '  +2.67%  '
$ws.Range("E49").Value = '  +2.55%  '
$ws.Range("D50").Value = '''0.07205'
$ws.Range("E50").Value = '  +0.33%  '
$ws.Range("D51").Value = '''80.61'
$ws.Range("E51").Value = '  +2.16%  '
